# Automatic update of files.
# Update the "Förändrad" (changed) date column (C) for all data rows (2-9)
# from serial date 45183 (2023-09-14) to 45184 (2023-09-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C9").Value = 45184
